$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021" column (P) to the indicator table, mirroring the
# formatting of the preceding "2020" column (O).

# Header cell (row 4): year label 2021, formatted like O4.
$ws.Range("P4").Value = 2021
$ws.Range("O4").Copy() | Out-Null
$ws.Range("P4").PasteSpecial(-4122) | Out-Null

# Data cell (row 5): participation rate value, formatted like O5.
$ws.Range("P5").Value = 80.9
$ws.Range("O5").Copy() | Out-Null
$ws.Range("P5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Match the saved selection state recorded in the sheet view.
$ws.Range("N10").Select() | Out-Null
